# Update the "Contact information" placeholder on the Thank-you slide
# (slide 10) with real contact details instead of the generic
# "F name L name" / "Title" / "Address" placeholders.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$sh = $s.Shapes.Item("Contact information")
$tr = $sh.TextFrame.TextRange

function Replace-SubText($range, [string]$search, [string]$replacement) {
    $current = $range.Text
    $idx = $current.IndexOf($search)
    if ($idx -lt 0) {
        return
    }
    $target = $range.Characters($idx + 1, $search.Length)
    $target.Text = $replacement
}

Replace-SubText $tr "F name L name" "Suresh Babu"
Replace-SubText $tr "Title" "Architect"
Replace-SubText $tr "Address" "Sap America"
